# Problem 50 Java Solution Backtracking 66.45% 91%
# Append a new performance data row (row 40) to the "performance" sheet,
# mirroring the existing rows: column A = runtime percentile (0.00% format),
# column B = memory percentile (0% format), column C = shared string label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: 66.45% runtime, 91% memory, "Backtracking" label.
$ws.Range("A40").Value = 0.6645
$ws.Range("A40").NumberFormat = "0.00%"

$ws.Range("B40").Value = 0.91
$ws.Range("B40").NumberFormat = "0%"

$ws.Range("C40").Value = "Backtracking"

# Reflect the scrolled/selected view state seen after adding the row.
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("D40").Select()
